$wb = $excel.ActiveWorkbook

# Row 17 on ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 483.5303
$ws.Range("J17").Value = 343.32758
$ws.Range("L17").Value = 1029.98274
$ws.Range("N17").Value = -1365.98274

# Row 28 on ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 739.8571
$ws.Range("I28").Value = 900.0909
$ws.Range("J28").Value = 563.6
$ws.Range("K28").Value = 900.0909
$ws.Range("L28").Value = 563.6
$ws.Range("M28").Value = -415.0909
$ws.Range("N28").Value = -1533.6

# Row 64 on ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4000
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()

# Row 67 on ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 4000
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()

# Row 76 on ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3200
$ws.Range("I76").Value = 3200
$ws.Range("K76").Value = 3200
$ws.Range("M76").Value = -2885

# Row 79 on ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 3200
$ws.Range("I79").Value = 3200
$ws.Range("K79").Value = 3200
$ws.Range("M79").Value = -2108

# Row 93 on ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H93").Value = 33523.69
$ws.Range("J93").Value = 33523.69
$ws.Range("L93").Value = 33523.69
$ws.Range("N93").Value = -38515.69

# Row 103 on ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 8101.0835
$ws.Range("I103").Value = 1560.8
$ws.Range("J103").Value = 12772.714
$ws.Range("K103").Value = 4682.4
$ws.Range("L103").Value = 38318.142
$ws.Range("M103").Value = -4096.4
$ws.Range("N103").Value = -39490.142

# Row 107 on ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 1484.55
$ws.Range("I107").Value = 1998.6364
$ws.Range("J107").Value = 856.2222
$ws.Range("K107").Value = 1998.6364
$ws.Range("L107").Value = 856.2222
$ws.Range("M107").Value = -78.6364000000001
$ws.Range("N107").Value = -4696.2222

# Row 129 on ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 866.81
$ws.Range("J129").Value = 898.4787
$ws.Range("L129").Value = 2695.4361
$ws.Range("N129").Value = -12695.4361

# Row 138 on ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 5351.033
$ws.Range("I138").Value = 718.9231
$ws.Range("J138").Value = 7203.877
$ws.Range("K138").Value = 2156.7693
$ws.Range("L138").Value = 21611.631
$ws.Range("M138").Value = 2983.2307
$ws.Range("N138").Value = -31891.631

# Row 141 on ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 6862.472
$ws.Range("I141").Value = 7564.8667
$ws.Range("J141").Value = 3350.5
$ws.Range("K141").Value = 22694.6001
$ws.Range("L141").Value = 10051.5
$ws.Range("M141").Value = -17514.6001
$ws.Range("N141").Value = -20411.5

# Row 63 on ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 11547018
$ws.Range("I63").Value = 27704442
$ws.Range("J63").Value = 5999.5713
$ws.Range("K63").Value = 27704442
$ws.Range("L63").Value = 5999.5713
$ws.Range("M63").Value = -27703756
$ws.Range("N63").Value = -7371.5713

# Row 66 on ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 11547018
$ws.Range("I66").Value = 27704442
$ws.Range("J66").Value = 5999.5713
$ws.Range("K66").Value = 138522210
$ws.Range("L66").Value = 29997.8565
$ws.Range("M66").Value = -138518778
$ws.Range("N66").Value = -36861.85649999999

# Row 74 on ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 6052.95
$ws.Range("I74").Value = 7066.231
$ws.Range("J74").Value = 4171.143
$ws.Range("K74").Value = 7066.231
$ws.Range("L74").Value = 4171.143
$ws.Range("M74").Value = -6192.231
$ws.Range("N74").Value = -5919.143

# Row 77 on ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 6052.95
$ws.Range("I77").Value = 7066.231
$ws.Range("J77").Value = 4171.143
$ws.Range("K77").Value = 35331.155
$ws.Range("L77").Value = 20855.715
$ws.Range("M77").Value = -30963.155
$ws.Range("N77").Value = -29591.715

# Row 80 on ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H80").Value = 32326.727
$ws.Range("J80").Value = 32326.727
$ws.Range("L80").Value = 32326.727
$ws.Range("N80").Value = -34322.727

# Row 83 on ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H83").Value = 32326.727
$ws.Range("J83").Value = 32326.727
$ws.Range("L83").Value = 96980.181
$ws.Range("N83").Value = -106964.181

# Row 122 on ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2286.0625
$ws.Range("I122").Value = 1123.909
$ws.Range("J122").Value = 4842.8
$ws.Range("K122").Value = 3371.727
$ws.Range("L122").Value = 14528.4
$ws.Range("M122").Value = -921.7270000000003
$ws.Range("N122").Value = -19428.4

# Row 132 on ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2718.2942
$ws.Range("I132").Value = 1351.0834
$ws.Range("K132").Value = 4053.2502
$ws.Range("M132").Value = -1523.2502

# Row 137 on ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H137").Value = 44544
$ws.Range("J137").Value = 44544
$ws.Range("L137").Value = 44544
$ws.Range("N137").Value = -54744

# Row 139 on ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 43191.54
$ws.Range("J139").Value = 43191.54
$ws.Range("L139").Value = 43191.54
$ws.Range("N139").Value = -53471.54

# Row 94 on BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2277.5
$ws.Range("I94").Value = 2050
$ws.Range("J94").Value = 2505
$ws.Range("K94").Value = 2050
$ws.Range("L94").Value = 2505
$ws.Range("M94").Value = -1599
$ws.Range("N94").Value = -3407

# Row 135 on BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H135").Value = 41330
$ws.Range("J135").Value = 41330
$ws.Range("L135").Value = 41330
$ws.Range("N135").Value = -51470

# Row 137 on BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H137").Value = 45500
$ws.Range("J137").Value = 45500
$ws.Range("L137").Value = 45500
$ws.Range("N137").Value = -55700

# Row 138 on BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H138").Value = 41353.2
$ws.Range("J138").Value = 41353.2
$ws.Range("L138").Value = 41353.2
$ws.Range("N138").Value = -51633.2

# Row 105 on CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1876.4762
$ws.Range("I105").Value = 1916.6154
$ws.Range("J105").Value = 1811.25
$ws.Range("K105").Value = 1916.6154
$ws.Range("L105").Value = 1811.25
$ws.Range("M105").Value = -169.6153999999999
$ws.Range("N105").Value = -5305.25

# Row 5 on CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 636339.5
$ws.Range("J5").Value = 1113226
$ws.Range("L5").Value = 3339678
$ws.Range("N5").Value = -3339902

# Row 134 on CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 3405.6216
$ws.Range("I134").Value = 2696.32
$ws.Range("K134").Value = 8088.960000000001
$ws.Range("M134").Value = -3018.960000000001

# Row 135 on CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 636339.5
$ws.Range("J135").Value = 1113226
$ws.Range("L135").Value = 10019034
$ws.Range("N135").Value = -10024104

# Row 80 on GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 50002384
$ws.Range("I80").Value = 83335370
$ws.Range("K80").Value = 83335370
$ws.Range("M80").Value = -83334372

# Row 83 on GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 50002384
$ws.Range("I83").Value = 83335370
$ws.Range("K83").Value = 416676850
$ws.Range("M83").Value = -416671858

# Row 97 on GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1492.3077
$ws.Range("I97").Value = 1491.8182
$ws.Range("J97").Value = 1495
$ws.Range("K97").Value = 1491.8182
$ws.Range("L97").Value = 1495
$ws.Range("M97").Value = -995.8181999999999
$ws.Range("N97").Value = -2487

# Row 126 on GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3323.26
$ws.Range("I126").Value = 2944.16
$ws.Range("J126").Value = 4460.56
$ws.Range("K126").Value = 8832.48
$ws.Range("L126").Value = 13381.68
$ws.Range("M126").Value = -6362.48
$ws.Range("N126").Value = -18321.68

# Row 137 on GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H137").Value = 39316.668
$ws.Range("J137").Value = 39316.668
$ws.Range("L137").Value = 39316.668
$ws.Range("N137").Value = -49516.668

# Row 2 on LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 1252251.5
$ws.Range("J2").Value = 1252251.5
$ws.Range("L2").Value = 1252251.5
$ws.Range("N2").Value = -1252475.5

# Row 100 on LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 3017.25
$ws.Range("I100").Value = 1343.2858
$ws.Range("K100").Value = 1343.2858
$ws.Range("M100").Value = -802.2858000000001

# Row 80 on WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 45960
$ws.Range("J80").Value = 45960
$ws.Range("L80").Value = 45960
$ws.Range("N80").Value = -47956

# Row 83 on WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H83").Value = 45960
$ws.Range("J83").Value = 45960
$ws.Range("L83").Value = 137880
$ws.Range("N83").Value = -147864

# Row 136 on WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3329.162
$ws.Range("I136").Value = 1974.96
$ws.Range("J136").Value = 6150.4165
$ws.Range("K136").Value = 5924.88
$ws.Range("L136").Value = 18451.2495
$ws.Range("M136").Value = -3374.88
$ws.Range("N136").Value = -23551.2495
